# LOQ4019.xlsx content update
# The original sheet had its row labels (column A) and long-form content
# (columns B/C) out of sync by a couple of rows starting at row 13. This
# script fixes the alignment: labels shift down a row starting at row 14,
# a new row 13 carries the professor's name, and the real long-form texts
# (objectives, syllabus, evaluation criteria, bibliography) are filled in
# where they belong. A new row 24 is added for the prerequisites text that
# used to sit in row 23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- long text values -------------------------------------------------
$objetivos = "1) Gerais: - Levar os estudantes a compreenderem os mecanismos de obtenção da influencia de diversos fatores ( variáveis independentes de um processo) sobre as variáveis resposta ( dependentes), através da análise multivariada.2) Específicos: - Ao final do curso os educandos devem:? Saber planejar e executar um experimento fatorial completo e fracionado? Saber analisar os resultados propondo a condição de melhor ajuste que otimiza os valores da variável resposta na região experimental estudada? Dominar, pelo menos, um software comercial sobre o assunto? Saber modelar um processo , com base em dados empíricos"
$programaResumido = "IntroduçãoExperimentação convencionalExperimentos Fatoriais completosExperimentos Fatoriais fracionadosAnálise de variânciaMetodologia de superfície de respostaMétodo de Taguchi"
$criterio = "serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso.A média da disciplina será a média aritmética das duas provas."
$bibliografia = "1) Planejamento e otimização de Experimentos. Roy E. Bruns, Edit. UNICAMP, 19963) Design and Analysis of Experiments, Douglas C. Montgomery, 6th edition, wiley, 20054) Designing for Quality  Robert H. Lochner  Ed. Quality Press, 19945) Statistics for Experimenter. Box & Hunter"

# ---- row 10: Objetivos/Objectives content ------------------------------
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# ---- row 13 loses its "Programa resumido:" label and instead carries ---
# ---- the professor's name (previously mis-placed on row 10) -----------
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5840535 - Messias Borges Silva"
$ws.Range("C13").Value = "5840535 - Messias Borges Silva"
$ws.Rows.Item(13).AutoFit()

# ---- row 14: "Programa resumido:" label + its real content ------------
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B20").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido
$ws.Rows.Item(14).RowHeight = 60

# ---- row 15: "Short syllabus:" label only ------------------------------
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
$ws.Rows.Item(15).RowHeight = 60

# ---- row 16: "Programa:" label + the same long syllabus content -------
$ws.Range("A16").Value = "Programa:"
$ws.Range("B20").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = $programaResumido
$ws.Range("C16").Value = $programaResumido
$ws.Rows.Item(16).RowHeight = 120

# ---- row 17: "Syllabus:" label only ------------------------------------
$ws.Range("A17").Value = "Syllabus:"
$ws.Rows.Item(17).RowHeight = 120

# ---- row 18: "Avaliação:" label only, no more mis-placed name ---------
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).AutoFit()

# ---- row 19: "Método:" label, keeps "2 provas escritas" ----------------
$ws.Range("A19").Value = "Método:"

# ---- row 20: "Critério:" label + the real grading-criteria text -------
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# ---- row 21: "Norma de recuperação:" label, keeps its text ------------
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Rows.Item(21).RowHeight = 60

# ---- row 22: "Bibliografia:" label + the real bibliography text -------
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B20").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
$ws.Rows.Item(22).RowHeight = 120

# ---- row 24: the prerequisites text that used to be on row 23 ---------
$ws.Range("B20").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("C20").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("B24").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Rows.Item(24).RowHeight = 30

# ---- row 23: "Requisitos:" label only (text now lives on row 24) ------
$ws.Range("A23").Value = "Requisitos:"
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()
$ws.Rows.Item(23).AutoFit()
